# Generate Report for Archive
#
# This script updates the "Ready for handoff" status text to "In Translation"
# across all worksheets, and narrows the related "Status"/language columns
# that displayed that (now longer-obsolete) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the status cells wherever "Ready for handoff" appeared.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Narrow the columns that held the status text (17.216 -> 13.410 chars).
# (12.5 is the nearest character width the host's pixel-grid snapping
# resolves to 13.41 chars of on-disk column width.)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
